$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: introduce new shared strings in the exact order used by the
# original authoring tool, so the rebuilt shared-string table lines up with
# the target index numbering ---
$ws.Cells.Item(14, 1).Value = "SelectPage"
$ws.Cells.Item(14, 3).Value = "ControlGroupSelectView_AvailabilityInputSelectView_RadioButtonMkt1Fare1"
$ws.Cells.Item(19, 2).Value = "Plus"
$ws.Cells.Item(20, 2).Value = "Max"
$ws.Cells.Item(21, 2).Value = "Agree"
$ws.Cells.Item(22, 2).Value = "SelectButton"
$ws.Cells.Item(1, 4).Value  = "Locator"
$ws.Cells.Item(2, 4).Value  = "id"
$ws.Cells.Item(21, 3).Value = "summary-user-acknowledge-checkbox"
$ws.Cells.Item(22, 3).Value = "ControlGroupSelectView_ButtonSubmit"
$ws.Cells.Item(18, 4).Value = "cssSelector"
$ws.Cells.Item(19, 3).Value = "button[data-for=bundlePlus1]"
$ws.Cells.Item(20, 3).Value = "button[data-for=bundleMax1]"
$ws.Cells.Item(18, 2).Value = "Starter"
$ws.Cells.Item(18, 3).Value = "div[class=starterbtn-option]"
$ws.Cells.Item(15, 3).Value = "ControlGroupSelectView_AvailabilityInputSelectView_RadioButtonMkt1Fare4"
$ws.Cells.Item(16, 3).Value = "ControlGroupSelectView_AvailabilityInputSelectView_RadioButtonMkt1Fare7"
$ws.Cells.Item(17, 3).Value = "ControlGroupSelectView_AvailabilityInputSelectView_RadioButtonMkt1Fare10"
$ws.Cells.Item(17, 2).Value = "FlightChoice10"
$ws.Cells.Item(14, 2).Value = "FlightChoice1"
$ws.Cells.Item(15, 2).Value = "FlightChoice4"
$ws.Cells.Item(16, 2).Value = "FlightChoice7"

# --- Step 2: fill in the remaining cells that reuse strings already present
# in the shared-string table (column A "SelectPage" + Locator column "id" /
# "cssSelector" repeats) ---

# Existing SearchPage rows (2-13): Locator column -> "id"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Value = "id"
}

# SelectPage rows 14-22, column A
for ($r = 14; $r -le 22; $r++) {
    $ws.Cells.Item($r, 1).Value = "SelectPage"
}

# Remaining Locator values for rows 14-22
$ws.Cells.Item(14, 4).Value = "id"
$ws.Cells.Item(15, 4).Value = "id"
$ws.Cells.Item(16, 4).Value = "id"
$ws.Cells.Item(17, 4).Value = "id"
$ws.Cells.Item(19, 4).Value = "cssSelector"
$ws.Cells.Item(20, 4).Value = "cssSelector"
$ws.Cells.Item(21, 4).Value = "id"
$ws.Cells.Item(22, 4).Value = "id"

# Selection cell
$ws.Range("C11").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
